# Applies the "Add files via upload" revision to the Assignment & Timeline
# workbook: refreshed Start/End dates (and their derived Duration formulas)
# for several tasks, a new "未執行" status note on row 34, and the
# view/formatting tweaks (zoom, frozen-pane selection, row heights, column
# widths) that came along with the re-save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Date($addr, $y, $m, $d) {
    $ws.Range($addr).Value = (Get-Date -Year $y -Month $m -Day $d).Date
}

# --- Start/End date corrections (rows 21-28) --------------------------------
Set-Date "C21" 2016 7 22
Set-Date "D21" 2016 7 28

Set-Date "C22" 2016 8 13
Set-Date "D22" 2016 8 14

Set-Date "C23" 2016 8 15
Set-Date "D23" 2016 8 18

Set-Date "C24" 2016 8 19
Set-Date "D24" 2016 8 19

Set-Date "C25" 2016 8 19
Set-Date "D25" 2016 8 19

Set-Date "C26" 2016 8 19
Set-Date "D26" 2016 8 19

Set-Date "C28" 2016 8 26
Set-Date "D28" 2016 8 26

# --- Newly populated Start/End dates (rows 29-33) ---------------------------
Set-Date "C29" 2016 8 27
Set-Date "D29" 2016 8 28

Set-Date "C30" 2016 8 21
Set-Date "D30" 2016 8 21

Set-Date "C31" 2016 7 22
Set-Date "D31" 2016 9 8

Set-Date "C32" 2016 7 22
Set-Date "D32" 2016 9 8

Set-Date "C33" 2016 7 22
Set-Date "D33" 2016 9 8

# --- New status note -----------------------------------------------------
$ws.Range("B34").Value = "未執行"

# --- Row heights (auto-computed wrap heights shifted slightly) -------------
$ws.Rows.Item(1).RowHeight = 48.6
$ws.Rows.Item(9).RowHeight = 64.8

# --- Column widths (re-fit to the refreshed content) ------------------------
$ws.Columns.Item(1).ColumnWidth = 82
$ws.Columns.Item(2).ColumnWidth = 36.142857142857146
$ws.Columns.Item(3).ColumnWidth = 9.428571428571429
$ws.Columns.Item(4).ColumnWidth = 9.428571428571429
$ws.Columns.Item(5).ColumnWidth = 9.285714285714286
$ws.Columns.Item(6).ColumnWidth = 15.142857142857142
$ws.Columns.Item(7).ColumnWidth = 14

# --- View state: zoom, frozen pane, selection --------------------------------
$ws.Range("A32").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 96
$ws.Range("F46").Select()
